{"js": "// Apply the benchmark-results update to the single-column results table.\n// Each table row holds one cell whose body text is replaced (as a single\n// run) with the new value, preserving the existing run formatting\n// (Times New Roman, sz 22) that Word keeps on the paragraph's first run.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Map of 0-indexed row number -> new cell text.\nconst rowUpdates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"614\",\n  5: \"0.00077\",\n  6: \"0.00018\",\n  7: \"0.00006\",\n  8: \"0.00029\",\n  9: \"0.00038\",\n  10: \"0.00048\",\n  11: \"0.12641\",\n  43: \"99.95\",\n  44: \"0.13\",\n  45: \"238\",\n};\n\nfor (const rowIndexStr of Object.keys(rowUpdates)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const newText = rowUpdates[rowIndexStr];\n\n  const cell = table.getCell(rowIndex, 0);\n  const paragraph = cell.body.paragraphs.getFirst();\n  const range = paragraph.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply the benchmark-results update to the single-column results table.\n# For each affected row, replace the cell's Range.Text with the new value\n# (this preserves the existing run formatting - Times New Roman, sz 22 -\n# carried by the cell's paragraph mark, and collapses any tab-separated\n# runs in a row down to a single run of new text).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowUpdates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"614\"\n    6  = \"0.00077\"\n    7  = \"0.00018\"\n    8  = \"0.00006\"\n    9  = \"0.00029\"\n    10 = \"0.00038\"\n    11 = \"0.00048\"\n    12 = \"0.12641\"\n    44 = \"99.95\"\n    45 = \"0.13\"\n    46 = \"238\"\n}\n\nforeach ($rowIndex in $rowUpdates.Keys) {\n    $cell = $t.Cell($rowIndex, 1)\n    $cell.Range.Text = $rowUpdates[$rowIndex]\n}\n"}
